$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.353.20'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.12%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.093.30'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.57%  '

# Row 4
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.22'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.19%  '

# Row 6
$ws.Range("E6").Value = '  +0.64%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.31'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +20.10%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '62.37'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.98%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.379'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.42%  '

# Row 11
$ws.Range("E11").Value = '  +4.45%  '

# Row 12
$ws.Range("E12").Value = '  +7.91%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.32'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +5.92%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.400.45'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.66%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.838'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.92%  '

# Row 16
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.27'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +7.24%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.097.93'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.74%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.260.93'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.84%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.29'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.71%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.54'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +14.91%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0854'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +4.85%  '

# Row 22
$ws.Range("E22").Value = '  +1.28%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.27'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +7.62%  '

# Row 24
$ws.Range("E24").Value = '  -0.08%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.47'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.87%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.51'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.70%  '

# Row 27
$ws.Range("E27").Value = '  +4.47%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.95'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.83%  '

# Row 29
$ws.Range("E29").Value = '  +4.73%  '

# Row 30
$ws.Range("E30").Value = '  +2.77%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.59'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +10.25%  '

# Row 32
$ws.Range("E32").Value = '  +23.37%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.56'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +4.97%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0628'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +7.99%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0907'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.58%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.28'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +7.17%  '

# Row 37
$ws.Range("E37").Value = '  +0.08%  '

# Row 38
$ws.Range("E38").Value = '  -2.95%  '

# Row 39
$ws.Range("E39").Value = '  -3.38%  '

# Row 40
$ws.Range("E40").Value = '  +3.36%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.81'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +139.88%  '

# Row 42
$ws.Range("E42").Value = '  +6.60%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.98'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +13.19%  '

# Row 44
$ws.Range("E44").Value = '  +5.35%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.42'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.55%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0959'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +17.05%  '

# Row 47
$ws.Range("E47").Value = '  +0.58%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.335.64'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.42%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.94'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.78%  '

# Row 50
$ws.Range("E50").Value = '  +7.29%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.98'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +13.53%  '
